$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the styling of the existing header row (A1:E1) for the new header cells
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New header cells F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean outlier flag values for F2:H12
$values = @(
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
